$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue #394: remove 'mrsow' and 'tsland' (rows 51-52) from the
# pre-identified-missing-variables list, along with the blank spacer row
# just above them (row 50) and one blank filler row just below them
# (row 53), so that the "LImon" block that used to start at row 55 now
# directly follows the row-49 spacer at row 51, and everything further
# down shifts up accordingly.
$ws.Rows("50:53").EntireRow.Delete()

# Reflect the resulting selection/scroll position after the rows above
# were removed.
$ws.Range("A50").Select()
